# Weekly update: insert a new daily price record for "Arveja Verde" at
# Femacal de La Calera (Coquimbo) as row 8, pushing existing historical
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 8; existing rows 8..62 shift to 9..63
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the latest record
$ws.Range("A8").Value2 = 3
$ws.Range("B8").Value2 = "Femacal de La Calera"
$ws.Range("C8").Value2 = "Coquimbo"
$ws.Range("D8").Value2 = 44831
$ws.Range("E8").Value2 = 5
$ws.Range("F8").Value2 = 100112022
$ws.Range("G8").Value2 = "Arveja Verde"
$ws.Range("H8").Value2 = "Perfection"
$ws.Range("I8").Value2 = "Primera"
$ws.Range("J8").Value2 = 45
$ws.Range("K8").Value2 = 28000
$ws.Range("L8").Value2 = 28000
$ws.Range("M8").Value2 = 28000
$ws.Range("N8").Value2 = "$/saco 25 kilos"
$ws.Range("O8").Value2 = "Provincia de Limarí"
$ws.Range("P8").Value2 = 1120
$ws.Range("Q8").Value2 = 25
$ws.Range("R8").Value2 = "Hortaliza"

# Apply the same date-number-format style used by the other date cells in column D
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
